$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 8343
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 8343
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 25029
$ws.Range("N96").Value = -27775
$ws.Range("M96").ClearContents()

$ws.Range("H125").Value = 711.86957
$ws.Range("I125").Value = 474.07144
$ws.Range("J125").Value = 1081.7778
$ws.Range("K125").Value = 4266.64296
$ws.Range("L125").Value = 9736.0002
$ws.Range("M125").Value = -1806.64296
$ws.Range("N125").Value = -14656.0002

$ws.Range("H132").Value = 5207.75
$ws.Range("I132").Value = 4587.1055
$ws.Range("J132").Value = 17000
$ws.Range("K132").Value = 13761.3165
$ws.Range("L132").Value = 51000
$ws.Range("M132").Value = -11231.3165
$ws.Range("N132").Value = -56060

$ws.Range("H138").Value = 2541.5645
$ws.Range("I138").Value = 3206.9167
$ws.Range("J138").Value = 2381.88
$ws.Range("K138").Value = 9620.750100000001
$ws.Range("L138").Value = 7145.64
$ws.Range("M138").Value = -4480.750100000001
$ws.Range("N138").Value = -17425.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 374507.28
$ws.Range("I32").Value = 447842.53
$ws.Range("J32").Value = 51832.266
$ws.Range("K32").Value = 447842.53
$ws.Range("L32").Value = 51832.266
$ws.Range("M32").Value = -447555.53
$ws.Range("N32").Value = -52406.266

$ws.Range("H61").Value = 9526061
$ws.Range("I61").Value = 20834828
$ws.Range("J61").Value = 2888.3157
$ws.Range("K61").Value = 20834828
$ws.Range("L61").Value = 2888.3157
$ws.Range("M61").Value = -20834616
$ws.Range("N61").Value = -3312.3157

$ws.Range("H74").Value = 686.7838
$ws.Range("I74").Value = 425.85715
$ws.Range("J74").Value = 1029.25
$ws.Range("K74").Value = 425.85715
$ws.Range("L74").Value = 1029.25
$ws.Range("M74").Value = 448.14285
$ws.Range("N74").Value = -2777.25

$ws.Range("H77").Value = 686.7838
$ws.Range("I77").Value = 425.85715
$ws.Range("J77").Value = 1029.25
$ws.Range("K77").Value = 2129.28575
$ws.Range("L77").Value = 5146.25
$ws.Range("M77").Value = 2238.71425
$ws.Range("N77").Value = -13882.25

$ws.Range("H97").Value = 1173.2059
$ws.Range("I97").Value = 1167.12
$ws.Range("J97").Value = 1190.1111
$ws.Range("K97").Value = 1167.12
$ws.Range("L97").Value = 1190.1111
$ws.Range("M97").Value = -671.1199999999999

$ws.Range("H132").Value = 5671.1055
$ws.Range("I132").Value = 6365.615
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 19096.845
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -16566.845
$ws.Range("N132").Value = -17559.0005

$ws.Range("H133").Value = 70130.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70130.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70130.5
$ws.Range("N133").Value = -75190.5

$ws.Range("H136").Value = 9526061
$ws.Range("I136").Value = 20834828
$ws.Range("J136").Value = 2888.3157
$ws.Range("K136").Value = 62504484
$ws.Range("L136").Value = 8664.947100000001
$ws.Range("M136").Value = -62501934
$ws.Range("N136").Value = -13764.9471

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2296.95
$ws.Range("I134").Value = 2134.8
$ws.Range("J134").Value = 3432
$ws.Range("K134").Value = 6404.400000000001
$ws.Range("L134").Value = 10296
$ws.Range("M134").Value = -3869.400000000001
$ws.Range("N134").Value = -15366

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 30071
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30071
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30071
$ws.Range("N44").Value = -30955

$ws.Range("H58").Value = 2806
$ws.Range("I58").Value = 2806
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2806
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2603
$ws.Range("N58").ClearContents()

$ws.Range("H132").Value = 16670066
$ws.Range("I132").Value = 3098
$ws.Range("J132").Value = 33337034
$ws.Range("K132").Value = 9294
$ws.Range("L132").Value = 100011102
$ws.Range("M132").Value = -6764

$ws.Range("H136").Value = 2806
$ws.Range("I136").Value = 2806
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8418
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5868
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 13513957
$ws.Range("I34").Value = 241
$ws.Range("J34").Value = 16129515
$ws.Range("K34").Value = 723
$ws.Range("L34").Value = 48388545
$ws.Range("M34").Value = -639
$ws.Range("N34").Value = -48388713

$ws.Range("H48").Value = 7124.5
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7124.5
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 21373.5
$ws.Range("N48").Value = -21873.5
$ws.Range("M48").ClearContents()

$ws.Range("H55").Value = 2082.353
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2082.353
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 6247.059
$ws.Range("N55").Value = -6601.059

$ws.Range("H68").Value = 1388.5275
$ws.Range("I68").Value = 798.7273
$ws.Range("J68").Value = 1576.5797
$ws.Range("K68").Value = 2396.1819
$ws.Range("L68").Value = 4729.7391
$ws.Range("M68").Value = -1585.1819
$ws.Range("N68").Value = -6351.7391

$ws.Range("H71").Value = 1388.5275
$ws.Range("I71").Value = 798.7273
$ws.Range("J71").Value = 1576.5797
$ws.Range("K71").Value = 7188.545700000001
$ws.Range("L71").Value = 14189.2173
$ws.Range("M71").Value = -3132.545700000001
$ws.Range("N71").Value = -22301.2173

$ws.Range("H113").Value = 625.86487
$ws.Range("I113").Value = 397.4074
$ws.Range("J113").Value = 1242.7
$ws.Range("K113").Value = 1192.2222
$ws.Range("L113").Value = 3728.1
$ws.Range("M113").Value = 977.7778000000001
$ws.Range("N113").Value = -8068.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1515.4445
$ws.Range("I97").Value = 1663.4286
$ws.Range("J97").Value = 997.5
$ws.Range("K97").Value = 1663.4286
$ws.Range("L97").Value = 997.5
$ws.Range("M97").Value = -1167.4286
$ws.Range("N97").Value = -1989.5

$ws.Range("H132").Value = 2308.8708
$ws.Range("I132").Value = 1945.2693
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 5835.8079
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -3305.8079
$ws.Range("N132").Value = -17658.8

$ws.Range("H138").Value = 50429
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50429
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50429
$ws.Range("N138").Value = -60709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 333336130
$ws.Range("I7").Value = 500001700
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 500001700
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -500001588
$ws.Range("N7").Value = -5229

$ws.Range("H40").Value = 125003540
$ws.Range("I40").Value = 250001950
$ws.Range("J40").Value = 5125
$ws.Range("K40").Value = 250001950
$ws.Range("L40").Value = 5125
$ws.Range("M40").Value = -250001814

$ws.Range("H46").Value = 3028.8
$ws.Range("I46").Value = 894
$ws.Range("J46").Value = 3562.5
$ws.Range("K46").Value = 894
$ws.Range("L46").Value = 3562.5
$ws.Range("M46").Value = -706
$ws.Range("N46").Value = -3938.5

$ws.Range("H55").Value = 398.10715
$ws.Range("I55").Value = 143.38889
$ws.Range("J55").Value = 856.6
$ws.Range("K55").Value = 143.38889
$ws.Range("L55").Value = 856.6
$ws.Range("M55").Value = 29.61111
$ws.Range("N55").Value = -1202.6

$ws.Range("H61").Value = 3117.8333
$ws.Range("I61").Value = 2803.7273
$ws.Range("J61").Value = 3611.4285
$ws.Range("K61").Value = 2803.7273
$ws.Range("L61").Value = 3611.4285
$ws.Range("M61").Value = -2601.7273
$ws.Range("N61").Value = -4015.4285

$ws.Range("H113").Value = 3117.8333
$ws.Range("I113").Value = 2803.7273
$ws.Range("J113").Value = 3611.4285
$ws.Range("K113").Value = 2803.7273
$ws.Range("L113").Value = 3611.4285
$ws.Range("M113").Value = -633.7273
$ws.Range("N113").Value = -7951.4285

$ws.Range("H126").Value = 333336130
$ws.Range("I126").Value = 500001700
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 1500005100
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -1500002630
$ws.Range("N126").Value = -19955

$ws.Range("H132").Value = 5021.643
$ws.Range("I132").Value = 5554.6665
$ws.Range("J132").Value = 4621.875
$ws.Range("K132").Value = 16663.9995
$ws.Range("L132").Value = 13865.625
$ws.Range("M132").Value = -14133.9995
$ws.Range("N132").Value = -18925.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3099.9092
$ws.Range("I81").Value = 2910.1
$ws.Range("J81").Value = 4998
$ws.Range("K81").Value = 5820.2
$ws.Range("L81").Value = 9996
$ws.Range("M81").Value = -4759.2
$ws.Range("N81").Value = -12118

$ws.Range("H84").Value = 3099.9092
$ws.Range("I84").Value = 2910.1
$ws.Range("J84").Value = 4998
$ws.Range("K84").Value = 29101
$ws.Range("L84").Value = 49980
$ws.Range("M84").Value = -23797
$ws.Range("N84").Value = -60588

$ws.Range("H105").Value = 99307.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 99307.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 99307.5
$ws.Range("N105").Value = -106295.5

$ws.Range("H107").Value = 668.1539
$ws.Range("I107").Value = 496.14285
$ws.Range("J107").Value = 868.8333
$ws.Range("K107").Value = 1488.42855
$ws.Range("L107").Value = 2606.4999
$ws.Range("M107").Value = 431.5714499999999
$ws.Range("N107").Value = -6446.4999

$ws.Range("H114").Value = 56699
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 56699
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 56699
$ws.Range("N114").Value = -65377

$ws.Range("H117").Value = 98409
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 98409
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 98409
$ws.Range("N117").Value = -107587

$ws.Range("H126").Value = 1567
$ws.Range("I126").Value = 1178.8
$ws.Range("J126").Value = 2052.25
$ws.Range("K126").Value = 3536.4
$ws.Range("L126").Value = 6156.75
$ws.Range("M126").Value = -1066.4
$ws.Range("N126").Value = -11096.75

$ws.Range("H132").Value = 10060076
$ws.Range("I132").Value = 2900.2307
$ws.Range("J132").Value = 18231530
$ws.Range("K132").Value = 8700.6921
$ws.Range("L132").Value = 54694590
$ws.Range("M132").Value = -6170.6921
$ws.Range("N132").Value = -54699650

$ws.Range("H136").Value = 3004.6099
$ws.Range("I136").Value = 2823.76
$ws.Range("J136").Value = 3287.1875
$ws.Range("K136").Value = 8471.280000000001
$ws.Range("L136").Value = 9861.5625
$ws.Range("M136").Value = -5921.280000000001
$ws.Range("N136").Value = -14961.5625
